$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# Locate the old marker text inside paragraph 1 (first run) dynamically so we
# don't depend on hard-coded character offsets.
$markerRange = $p.Range.Duplicate
$found = $markerRange.Find.Execute("**ID__AFFARS_5314_topic_5__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Anything between the end of the marker text and the end of the paragraph
# (excluding the trailing pilcrow) is the extra " " run that must be dropped
# entirely rather than merged into the remaining run.
$paraEnd = $p.Range.End
if ($paraEnd - 1 -gt $markerRange.End) {
    $trailingRange = $d.Range($markerRange.End, $paraEnd - 1)
    $trailingRange.Delete()
}

# Replace the marker run's text in place (keeps its original rPr formatting).
$markerRange.Text = "**ID__AFFARS_5314_407_3__ID**"

# Add the paragraph border (space-only border, no visible line) to paragraph 1,
# matching the border already used by the later paragraphs in the document.
$p.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p.Range.ParagraphFormat.Borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p.Range.ParagraphFormat.LeftIndent = 11.25
